$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated cryptos list values (price / 1h volume, plus a couple of row swaps) ---

# Cells whose new text reads as a plain number (e.g. "0.999", "54.10").
# These columns are stored as text in the workbook (to preserve values like
# "71.224.61" that are not valid numbers), so we force text format before
# assigning, then restore the default "Normal" style so no extra formatting lingers.
$numericLookingValues = @{
    "D4" = "0.999"
    "D5" = "606.61"
    "D6" = "202.79"
    "D10" = "0.651"
    "D11" = "54.10"
    "D12" = "0.0000308"
    "D15" = "678.50"
    "D17" = "12.94"
    "D19" = "19.10"
    "D22" = "18.96"
    "D23" = "5.42"
    "D24" = "105.36"
    "D25" = "4.64"
    "D27" = "10.65"
    "D29" = "34.37"
    "D30" = "4.59"
    "D32" = "12.29"
    "D33" = "0.116"
    "D34" = "63.54"
    "D37" = "1.00"
    "D38" = "520.73"
    "D41" = "0.393"
    "D42" = "36.79"
    "D43" = "0.138"
    "D44" = "3.07"
    "D45" = "0.0460"
    "D46" = "3.47"
    "D48" = "8.70"
}

foreach ($addr in $numericLookingValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingValues[$addr]
    $cell.Style = "Normal"
}

# Remaining cells (coin names, links, percentage strings, multi-dot prices)
# already round-trip as text without any special handling.
$plainTextValues = @{
    "D2" = "71.224.61"
    "E2" = "  +2.52%  "
    "D3" = "3.645.48"
    "E3" = "  +4.15%  "
    "E4" = "  +0.13%  "
    "E5" = "  +1.21%  "
    "E6" = "  +4.39%  "
    "E7" = "  +1.49%  "
    "E8" = "  +0.05%  "
    "E9" = "  +10.37%  "
    "E10" = "  +1.15%  "
    "E11" = "  +2.21%  "
    "E12" = "  +3.30%  "
    "E13" = "  +2.28%  "
    "D14" = "4.226.07"
    "E14" = "  +4.25%  "
    "E15" = "  +13.64%  "
    "D16" = "71.278.86"
    "E16" = "  +2.43%  "
    "E17" = "  +2.80%  "
    "D18" = "3.652.11"
    "E18" = "  +4.11%  "
    "E19" = "  +1.05%  "
    "E20" = "  +0.38%  "
    "E21" = "  +2.45%  "
    "E22" = "  +6.41%  "
    "E23" = "  +2.28%  "
    "E24" = "  +2.85%  "
    "E25" = "  +1.19%  "
    "E26" = "  -1.44%  "
    "E27" = "  -1.03%  "
    "E28" = "  +5.51%  "
    "E29" = "  +4.48%  "
    "E30" = "  +7.32%  "
    "E31" = "  +5.00%  "
    "E32" = "  +0.28%  "
    "E33" = "  +2.12%  "
    "E34" = "  +0.68%  "
    "D35" = "0.0₃0879"
    "E35" = "  +8.89%  "
    "D36" = "3.927.07"
    "E36" = "  +4.79%  "
    "B37" = "Dai"
    "C37" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "E37" = "  +0.04%  "
    "B38" = "Bittensor"
    "C38" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "E38" = "  +5.71%  "
    "E39" = "  -4.18%  "
    "E40" = "  -0.55%  "
    "E41" = "  +1.35%  "
    "E42" = "  +2.59%  "
    "E43" = "  +4.36%  "
    "B44" = "ThetaToken"
    "C44" = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
    "E44" = "  +9.64%  "
    "B45" = "VeChain"
    "C45" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "E45" = "  +2.80%  "
    "E46" = "  +6.95%  "
    "E47" = "  +1.47%  "
    "E48" = "  +3.79%  "
    "E50" = "  +2.47%  "
    "E51" = "  +4.70%  "
}

foreach ($addr in $plainTextValues.Keys) {
    $ws.Range($addr).Value = $plainTextValues[$addr]
}

